# Fixed LEAC_plot_iter.py len(rate_table)-1 error.
# Correct the 2026 "First 55000" rate from 0.22 to 0.23 and restore the
# view to the top of the sheet with B2 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rates")

$ws.Range("B3").Value = 0.23

$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
